# Remove the standalone "No API key found..." sample-output paragraph
# that followed the Python code block (the <pStyle SourceCode>/<rStyle
# VerbatimChar> paragraph), while leaving the print(...) statement that
# produces that text untouched.

$d = $word.ActiveDocument
$target = "No API key found. Please set the OPENWEATHER_API_KEY environment variable."

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd("`r", "`a")
    if ($text -eq $target) {
        # Delete the whole paragraph, including its paragraph mark.
        $p.Range.Delete()
        break
    }
}
